$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188, shifting rows 188:249 down to 189:250
$ws.Rows("188:188").Insert()

# Populate the newly inserted row 188 with data
$ws.Range("A188").Value = 6
$ws.Range("B188").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C188").Value = "Metropolitana"
$ws.Range("D188").Value = 44726
$ws.Range("E188").Value = 13
$ws.Range("F188").Value = 100112026
$ws.Range("G188").Value = "Haba"
$ws.Range("H188").Value = "Sin especificar"
$ws.Range("I188").Value = "Primera"
$ws.Range("J188").Value = 250
$ws.Range("K188").Value = 18000
$ws.Range("L188").Value = 20000
$ws.Range("M188").Value = 18800
$ws.Range("N188").Value = '$/saco 25 kilos'
$ws.Range("O188").Value = "Provincia de Huasco"
$ws.Range("P188").Value = 752
$ws.Range("Q188").Value = 25
$ws.Range("R188").Value = "Hortaliza"
